# Apply edits described by the commit:
#   "Changed variance scaling off for PLSR. removed the algorithm for mean
#    centering as it is done when instantiating PLSR from sklearn."
#
# Concretely, in the workbook this means the calibration_stats_turbid sheet
# was regenerated/re-sorted (descending by RMSEC, column I), an AutoFilter
# was applied to it (matching the existing one already present on
# calibration_stats_sn), and the active sheet/selection bookkeeping moved
# from calibration_stats_sn to calibration_stats_turbid.

$wb = $excel.ActiveWorkbook

$wsTurbid = $wb.Worksheets.Item("calibration_stats_turbid")
$wsSn     = $wb.Worksheets.Item("calibration_stats_sn")

# --- Re-sort calibration_stats_turbid data rows by RMSEC (column I), descending ---
$dataRange = $wsTurbid.Range("A1:K37")
$dataRange.Sort($wsTurbid.Range("I1"), 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, $null, $null)

# --- Apply an AutoFilter to the header row (mirrors calibration_stats_sn) ---
$headerRange = $wsTurbid.Range("A1:K1")
$headerRange.AutoFilter()

# Hidden defined name Excel keeps alongside an AutoFilter'ed range.
$wsTurbid.Names.Add("_xlnm._FilterDatabase", "=calibration_stats_turbid!`$A`$1:`$K`$1", $false) | Out-Null

# --- Make calibration_stats_turbid the active sheet / tab, with M10 selected ---
$wsTurbid.Activate()
$wsTurbid.Range("M10").Select()

# --- calibration_stats_sn keeps its existing K2 selection but is no longer the active tab ---
$wsSn.Range("K2").Select()

# Re-activate calibration_stats_turbid so it ends up as the active/visible sheet
$wsTurbid.Activate()
